# Updates the recalculated extrapolation-calibration outputs (columns D-H)
# for the rows affected by removing the < USD 5 price quote from the
# calibration inputs (treated as noise).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> updated values for columns D (ABSM1_RN), E (M1_RN), F (CM2_RN),
# G (CMN3_RN), H (CMN4_RN)
$updates = @{
    3  = @{ D = 112558.6169966177;  E = -0.01793792615951216; F = 0.168228484785362;  G = -0.9858338634566688; H = 8.814173138256782 }
    7  = @{ D = 115377.222283621;   E = -0.02457030015719104; F = 0.205356130593579;   G = -0.5979058647301055; H = 5.50870602778383 }
    8  = @{ D = 115642.2723198344; E = -0.03782191976319748; F = 0.1984381254595149;  G = -1.281155195030867;  H = 8.772604062186625 }
    9  = @{ D = 117156.835540117;   E = -0.06728746462106372; F = 0.3233693267818493;  G = -1.844944676513989;  H = 11.36404209307687 }
    10 = @{ D = 118668.331208669;   E = -0.1046918920944829;  F = 0.4135241985041633;  G = -1.872670482640333;  H = 9.58161668687984 }
    13 = @{ D = 111755.7589752862; E = -0.02770118955574321; F = 0.1400204828551337;  G = -0.263608381758998;  H = 4.878534593564482 }
    16 = @{ D = 111666.5225405059; E = -0.05058030847475746; F = 0.164508674736923;   G = -0.6950933036113355; H = 5.012956620614996 }
    18 = @{ D = 111674.572228204;   E = -0.03193426404774774; F = 0.1520460463169616;  G = -0.7151408038088994; H = 4.665645595731196 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $vals.D
    $ws.Cells.Item($row, 5).Value = $vals.E
    $ws.Cells.Item($row, 6).Value = $vals.F
    $ws.Cells.Item($row, 7).Value = $vals.G
    $ws.Cells.Item($row, 8).Value = $vals.H
}
